$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J, shifting Pitches..Protection one column to the right
$ws.Range("J:J").Insert()

# Set the new header cell for the inserted column
$ws.Range("J1").Value = "Aid_Grade"

# Update the active selection to J3 (as in the saved file)
$ws.Range("J3").Select()
